# Update cryptos list snapshot (prices / 1h volume deltas) and fix the
# swapped RocketPoolETH / ordi rows, per the scheduled GitHub Actions
# refresh.
#
# Note: several "Price" values are plain decimal-looking strings
# (e.g. "74.56"). Excel's COM Range.Value setter auto-coerces those to
# numbers, but the source data stores them as text, so a leading
# apostrophe is used to force text entry for anything that would
# otherwise parse as a number; Excel strips the apostrophe from the
# stored value itself.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value  = "42.645.94"
$ws.Range("E2").Value  = "  -0.56%  "

$ws.Range("D3").Value  = "2.530.95"
$ws.Range("E3").Value  = "  -0.68%  "

$ws.Range("E4").Value  = "  +0.01%  "

$ws.Range("D5").Value  = "'315.82"
$ws.Range("E5").Value  = "  +4.02%  "

$ws.Range("D6").Value  = "'95.31"
$ws.Range("E6").Value  = "  -2.68%  "

$ws.Range("E7").Value  = "  +0.51%  "

$ws.Range("E8").Value  = "  +0.01%  "

$ws.Range("E9").Value  = "  -1.33%  "

$ws.Range("D10").Value = "'36.42"
$ws.Range("E10").Value = "  -1.59%  "

$ws.Range("D11").Value = "'0.0813"
$ws.Range("E11").Value = "  -1.53%  "

$ws.Range("D12").Value = "'7.73"
$ws.Range("E12").Value = "  -0.42%  "

$ws.Range("E13").Value = "  -2.77%  "

$ws.Range("D14").Value = "2.917.78"
$ws.Range("E14").Value = "  -0.80%  "

$ws.Range("D15").Value = "'15.66"
$ws.Range("E15").Value = "  +3.50%  "

$ws.Range("D16").Value = "2.555.18"
$ws.Range("E16").Value = "  -0.30%  "

$ws.Range("D17").Value = "'0.860"
$ws.Range("E17").Value = "  -1.76%  "

$ws.Range("D18").Value = "42.680.72"
$ws.Range("E18").Value = "  -0.54%  "

$ws.Range("D19").Value = "'13.05"
$ws.Range("E19").Value = "  -5.58%  "

$ws.Range("D20").Value = "'6.66"
$ws.Range("E20").Value = "  +1.26%  "

$ws.Range("D21").Value = "0.0₃0970"
$ws.Range("E21").Value = "  -2.41%  "

$ws.Range("D22").Value = "'71.34"
$ws.Range("E22").Value = "  -0.75%  "

$ws.Range("D23").Value = "'254.55"
$ws.Range("E23").Value = "  +0.05%  "

$ws.Range("D24").Value = "'2.98"
$ws.Range("E24").Value = "  +0.63%  "

$ws.Range("E25").Value = "  -1.69%  "

$ws.Range("D26").Value = "'27.60"
$ws.Range("E26").Value = "  -1.46%  "

$ws.Range("D27").Value = "'0.995"
$ws.Range("E27").Value = "  -0.49%  "

$ws.Range("D28").Value = "'2.32"
$ws.Range("E28").Value = "  +11.18%  "

$ws.Range("D29").Value = "'39.50"
$ws.Range("E29").Value = "  +4.68%  "

$ws.Range("D30").Value = "'10.08"
$ws.Range("E30").Value = "  -1.63%  "

$ws.Range("D31").Value = "'5.91"
$ws.Range("E31").Value = "  -3.87%  "

$ws.Range("D32").Value = "'156.39"
$ws.Range("E32").Value = "  -1.28%  "

$ws.Range("D33").Value = "'20.02"
$ws.Range("E33").Value = "  +2.11%  "

$ws.Range("D34").Value = "'3.36"
$ws.Range("E34").Value = "  +1.62%  "

$ws.Range("D35").Value = "'2.12"
$ws.Range("E35").Value = "  -1.40%  "

$ws.Range("D36").Value = "'0.0787"
$ws.Range("E36").Value = "  -1.91%  "

$ws.Range("E37").Value = "  -0.75%  "

$ws.Range("E38").Value = "  -3.08%  "

$ws.Range("D39").Value = "'24.46"
$ws.Range("E39").Value = "  -4.00%  "

$ws.Range("E40").Value = "  -0.02%  "

$ws.Range("D41").Value = "'2.20"
$ws.Range("E41").Value = "  +4.05%  "

$ws.Range("D42").Value = "'3.40"
$ws.Range("E42").Value = "  -0.53%  "

$ws.Range("E43").Value = "  -1.07%  "

$ws.Range("D44").Value = "'0.0304"
$ws.Range("E44").Value = "  -0.65%  "

$ws.Range("E45").Value = "  +0.02%  "

$ws.Range("D46").Value = "2.049.22"
$ws.Range("E46").Value = "  -2.15%  "

$ws.Range("D47").Value = "'86.03"
$ws.Range("E47").Value = "  -0.58%  "

$ws.Range("D48").Value = "'8.85"
$ws.Range("E48").Value = "  -1.24%  "

# Rows 49/50 swap places (RocketPoolETH <-> ordi) with refreshed figures;
# the rank column (A) is untouched.
$ws.Range("B49").Value = "ordi"
$ws.Range("C49").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D49").Value = "'74.56"
$ws.Range("E49").Value = "  -0.98%  "

$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").Value = "2.769.75"
$ws.Range("E50").Value = "  -1.00%  "

$ws.Range("D51").Value = "'0.191"
$ws.Range("E51").Value = "  -0.29%  "
